$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new schedule entry as row 50 (A: period start date, B: period end
# date, C: description) — continuing the existing table pattern.
$ws.Range("A50").Value = "2026/3/6"
$ws.Range("B50").Value = "2026/5/1"
$ws.Range("C50").Value = "第100期 第二代星途"

# Match the author's new viewport/selection state: scrolled down to row 36,
# with C51 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C51").Select()
